$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting Joy..I3 one column to the right
$ws.Columns("C:C").Insert()

# Set the new header cell's value
$ws.Range("C1").Value = "Neutral"

# Update the selected cell to match the target state
$ws.Range("E5").Select()
